# Updates the cryptos price/volume(1h) figures (and re-sorts a few rows
# whose relative ranking changed) to match the refreshed coinranking.com
# snapshot from the GitHub Actions run.
#
# Every assignment is written as a leading-apostrophe literal (the COM/VBA
# idiom for "store as text even though it parses as a number") because the
# Price column holds values such as "26.661.48" / "0.790" that must stay
# plain text - Excel would otherwise silently coerce them to numbers (e.g.
# 211.30 -> 211.3). The immediate `.Style = 'Normal'` afterwards strips the
# implicit "quote prefix" cell style that leading-apostrophe entry adds, so
# the cell keeps its original (default) style, exactly as in the workbook
# before the edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''26.661.48'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -0.57%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '''1.597.17'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -0.67%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '''  +0.13%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''211.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  +0.30%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('E6').Value = '''  +0.68%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = '''  +0.12%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = '''  -0.29%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = '''0.246'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  -1.23%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = '''19.69'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  -0.27%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''0.0838'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  -0.09%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = '''1.821.06'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -0.64%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = '''1.571.89'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  -2.11%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('E14').Value = '''  -0.68%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('E15').Value = '''  -1.85%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = '''64.96'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  +1.98%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = '''26.646.03'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  -0.50%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''0.0₃0728'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  -0.23%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '''209.64'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '''  -0.17%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = '''  +0.04%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''6.78'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  +0.36%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = '''  -0.39%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = '''  -1.54%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('E24').Value = '''  +0.54%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''146.45'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''  -0.06%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = '''  +0.16%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = '''  -4.07%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').Value = '''  +2.05%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''15.29'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  -0.35%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''0.0504'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  +0.41%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = '''  +0.33%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = '''  -1.02%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''0.664'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  -1.65%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = '''  -1.22%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = '''1.298.45'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  -1.47%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = '''2.45'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  +0.58%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = '''  -3.03%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('E38').Value = '''  -1.22%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = '''  +2.46%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').Value = '''  +0.10%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('B41').Value = '''TrustWalletToken'
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').Value = '''https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').Value = '''0.790'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  +0.06%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('B42').Value = '''FraxShare'
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').Value = '''https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').Value = '''5.37'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  +1.81%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('B43').Value = '''MXToken'
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = '''https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = '''2.20'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  +0.23%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''63.84'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  +1.38%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = '''1.733.99'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  -0.63%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('B46').Value = '''WEMIXToken'
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = '''https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = '''0.883'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  +7.99%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('B47').Value = '''Quant'
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = '''https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = '''90.07'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  +1.12%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''1.63'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E49').Value = '''  +2.18%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = '''  -1.26%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''7.49'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  +0.04%  '
$ws.Range('E51').Style = 'Normal'
